$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 30017
$ws.Range("B2").Value = "Pedro Miguel Cavalcanti"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Viagem de negocios"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45081
$ws.Range("G2").Value = 9831.379999999999

# Row 3
$ws.Range("A3").Value = 92839
$ws.Range("B3").Value = "Sr. Dom Sousa"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45078
$ws.Range("G3").Value = 9933.610000000001

# Row 4
$ws.Range("A4").Value = 45948
$ws.Range("B4").Value = "João Gabriel Abreu"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Viagem de negocios"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45099
$ws.Range("G4").Value = 2413.9

# Row 5
$ws.Range("A5").Value = 83601
$ws.Range("B5").Value = "Thiago Siqueira"
$ws.Range("C5").Value = "Juridico"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45103
$ws.Range("G5").Value = 6873.86

# Row 6
$ws.Range("A6").Value = 51735
$ws.Range("B6").Value = "Marina Borges"
$ws.Range("C6").Value = "Operacoes"
$ws.Range("D6").Value = "Consulta medica"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45103
$ws.Range("G6").Value = 6774.14

# Row 7
$ws.Range("A7").Value = 85022
$ws.Range("B7").Value = "Isabella Dias"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Viagem de negocios"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45095
$ws.Range("G7").Value = 4479.66

# Row 8
$ws.Range("A8").Value = 71972
$ws.Range("B8").Value = "Vicente Costa"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45095
$ws.Range("G8").Value = 7700.09

# Row 9
$ws.Range("A9").Value = 68267
$ws.Range("B9").Value = "Ana Clara Peixoto"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45103
$ws.Range("G9").Value = 9849.82

# Row 10
$ws.Range("A10").Value = 31164
$ws.Range("B10").Value = "Otto da Costa"
$ws.Range("C10").Value = "Juridico"
$ws.Range("D10").Value = "Viagem de negocios"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 45079
$ws.Range("G10").Value = 7309.7

# Row 11
$ws.Range("A11").Value = 94715
$ws.Range("B11").Value = "Bella Peixoto"
$ws.Range("C11").Value = "Recursos Humanos"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45100
$ws.Range("G11").Value = 5711.05
